# Updated symbol list on Sat Dec 17 22:52:57 UTC 2022 with GitHub Actions
#
# Applies the refreshed "Price" (col D) quotes, plus a data-fix that swaps
# the CEJI / BKEXToken rows (42 <-> 43) back into coinranking's sort order,
# including their Link (col C) and rank-key (col E) columns.
#
# All of these columns are stored as literal TEXT in the workbook (not
# numbers), so every numeric-looking value is written through a
# Text-number-format round trip ('@' -> Value -> restore Normal style) to
# stop Excel's automatic "this looks like a number" reinterpretation from
# flipping the cell to the Number type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the literal string into the cell without leaving the cell
    # tagged as Text-formatted afterwards (matches the source file, where
    # these cells carry no explicit number format / style).
    $cell = $ws.Range($range)
    $cell.NumberFormat = '@'
    $cell.Value = $value
    $cell.Style = 'Normal'
}

# Col D "Price" refresh (A1 ref -> new text value)
$prices = @{
    'D2'  = '238.97'
    'D3'  = '21.63'
    'D4'  = '5.367'
    'D5'  = '0.05561'
    'D6'  = '3.366'
    'D7'  = '6.457'
    'D8'  = '0.8052'
    'D9'  = '1.050'
    'D10' = '0.1400'
    'D11' = '0.07300'
    'D12' = '0.03265'
    'D13' = '0.02956'
    'D14' = '0.09239'
    'D15' = '0.001648'
    'D16' = '3.250'
    'D17' = '0.04759'
    'D18' = '0.0005705'
    'D19' = '0.006255'
    'D20' = '0.001048'
    'D21' = '0.003781'
    'D22' = '0.0001497'
    'D23' = '0.0004177'
    'D24' = '3.964'
    'D25' = '2.205'
    'D27' = '0.1294'
    'D40' = '0.04162'
    'D41' = '0.006996'
    'D42' = '0.1039'
    'D43' = '0.002903'
    'D44' = '0.008841'
    'D45' = '0.00005436'
    'D47' = '0.6793'
    'D48' = '0.03135'
    'D49' = '0.00002098'
}

foreach ($ref in $prices.Keys) {
    Set-TextValue $ref $prices[$ref]
}

# Row 42 / 43 swap: CEJI <-> BKEXToken (Coin, Link, rank-key columns)
$ws.Range('B42').Value = 'BKEXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('E42').Value = '41BKEXTokenBKK'

$ws.Range('B43').Value = 'CEJI'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range('E43').Value = '42CEJICEJI'
